$wb = $excel.ActiveWorkbook

# --- Bus sheet: update E5 value and selection ---
$busSheet = $wb.Worksheets.Item("Bus")
$busSheet.Range("E5").Value = 1
$busSheet.Activate()
$busSheet.Range("E5").Select()

# --- Device sheet: update formulas in C5 and D5, and selection ---
$deviceSheet = $wb.Worksheets.Item("Device")
$deviceSheet.Range("C5").Formula = "=1/(2*2*3.14*100*2*3.14)*2*3.14*60/2"
$deviceSheet.Range("D5").Formula = "=1/(2*2*3.14)*2*3.14*60"
$deviceSheet.Activate()
$deviceSheet.Range("D5").Select()
